$wb = $excel.ActiveWorkbook

# Rename the "temperature_c" sheet to "temperature"
$ws = $wb.Worksheets.Item("temperature_c")
$ws.Name = "temperature"

# Switch the active/selected sheet from "genotype" to "temperature"
$ws.Activate()
